$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some columns (L, and M on the newly-filled row 338) are formatted
# as Text (numFmtId 49, format code "@"). Writing a numeric literal into such
# a cell via COM stores it as a text value (mirroring real Excel behaviour),
# which is not what the source data needs - the underlying values must stay
# numeric so the dependent SUM/IF formulas keep working. To store a genuine
# number while leaving the cell's displayed format untouched, temporarily
# switch the format to a plain numeric one, set the value, then restore the
# original format string.
function Set-NumericValue($address, $value) {
    $rng = $ws.Range($address)
    $originalFormat = $rng.NumberFormat
    if ($originalFormat -eq "@") {
        $rng.NumberFormat = "0"
        $rng.Value = $value
        $rng.NumberFormat = $originalFormat
    } else {
        $rng.Value = $value
    }
}

# --- Updated daily figures for rows 300-337 (new deaths reported / new case counts) ---

# Row 300: one more hospital death reported
Set-NumericValue "L300" 3

# Row 307: one more hospital death reported
Set-NumericValue "L307" 2

# Row 319: one more hospital death reported
Set-NumericValue "L319" 4

# Row 325: new positive cases revised down by 1, one more hospital death reported
Set-NumericValue "C325" 119
Set-NumericValue "L325" 5

# Row 326: one more hospital death reported
Set-NumericValue "L326" 2

# Row 334: new positive cases revised up by 1
Set-NumericValue "C334" 71

# Row 335: new positive cases revised up by 2
Set-NumericValue "C335" 142

# Row 336: new positive cases revised up, one more hospital death reported
Set-NumericValue "C336" 86
Set-NumericValue "L336" 2

# Row 337: new positive cases revised up, one more hospital death reported
Set-NumericValue "C337" 57
Set-NumericValue "L337" 2

# --- Row 338 (2020-12-17) is populated for the first time: it used to be an
# empty trailing row with only the shared formulas and no input data. ---
Set-NumericValue "C338" 17
Set-NumericValue "E338" 16
Set-NumericValue "F338" 13
Set-NumericValue "G338" 110
Set-NumericValue "L338" 0
Set-NumericValue "M338" 0
